# Apply the edit described by the diff:
#  - Insert 3 new columns (City, Country, Domain) after "Station Name" (i.e. before the
#    old "Branches" column), shifting Branches..Non Tech skills from C..I to F..L.
#  - Rewrite all data rows (2-6) with corrected / additional values, including some rows
#    whose station id / order changed, and fix "missing"/"wrong" values (e.g. the
#    "Weekly Holidays" value is now lower-case "sunday" and many Tech/Non-Tech skill
#    values that used to be "#NA" / "# NA" are replaced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write an array of values into a single worksheet row starting at
# column A, using a true 2D COM-safe array so Excel treats it as one bulk
# write (keeps default cell styling).
# ---------------------------------------------------------------------------
function Set-RowValues {
    param(
        $Sheet,
        [int]$RowNum,
        [object[]]$Values
    )
    $n = $Values.Length
    $data = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $data[0, $i] = $Values[$i]
    }
    $startCell = $Sheet.Cells.Item($RowNum, 1)
    $endCell = $Sheet.Cells.Item($RowNum, $n)
    $rng = $Sheet.Range($startCell, $endCell)
    $rng.Value = $data
}

# ---------------------------------------------------------------------------
# 1) Insert three new blank columns before the old column C ("Branches"),
#    shifting everything from C.. onward three places to the right.
# ---------------------------------------------------------------------------
$ws.Range("C:E").Insert()

# ---------------------------------------------------------------------------
# 2) Header row (row 1): fill in the three new column headers.
#    (Existing headers in F1:L1 already carry over correctly from the insert.)
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "City"
$ws.Cells.Item(1, 4).Value = "Country"
$ws.Cells.Item(1, 5).Value = "Domain"

# ---------------------------------------------------------------------------
# 3) Data rows: overwrite rows 2-6 completely with the corrected data set.
#    Columns: A Station ID | B Station Name | C City | D Country | E Domain |
#             F Branches | G Accomodation | H Timings | I Weekly Holidays |
#             J Stipend for Single Degree | K Tech skills | L Non Tech skills
# ---------------------------------------------------------------------------
Set-RowValues $ws 2 @(
    6865,
    "505 Army Base Workshop-CS/IT",
    "Delhi",
    "Delhi",
    "Mechanical",
    "A7",
    "NO",
    "9:00 AM to 2:00 PM",
    "sunday",
    "0",
    "Software Development",
    "Analytical Skill"
)

Set-RowValues $ws 3 @(
    6636,
    "505 Army Base Workshop-Electro Mechanics",
    "New Delhi",
    "Delhi",
    "Electronics",
    "A3, A4, A8, AA",
    "NO",
    "9:00 AM to 2:00 PM",
    "sunday",
    "0",
    "Electronics,Mechanical Processes,Mechanical & Electrical forces on Transformers,Mechanical,Knowledge of key Mechanical Equipment",
    "Analytical Skill,Excellent communication"
)

Set-RowValues $ws 4 @(
    6846,
    "505 Army Base Workshop-Manufacturing",
    "New Delhi",
    "Delhi",
    "Manufacturing",
    "A4, AB",
    "NO",
    "9:00 AM to 2:00 PM",
    "sunday",
    "0",
    "Mechanical Processes,Mechanical,Manufacturing,Manufacturing & Materials,Manufacturing Operation Tools",
    "Analytical Skill"
)

Set-RowValues $ws 5 @(
    6845,
    "505 Army Base Workshop-Mechanical",
    "New Delhi",
    "Delhi",
    "Electronics",
    "A4",
    "NO",
    "9:00 AM to 2:00 PM",
    "sunday",
    "0",
    "Mechanical,Mechanical & Electrical forces on Transformers",
    "Analytical Skill"
)

Set-RowValues $ws 6 @(
    6866,
    "505 Army Base Workshop-Mechatronics",
    "Delhi",
    "Delhi",
    "Mechanical",
    "A4, A7, A8, AA",
    "NO",
    "9:00 AM to 2:00 PM",
    "sunday",
    "0",
    "Mechanical,Mechanical & Electrical forces on Transformers,Mechanical Processes",
    "Analytical Skill"
)

# ---------------------------------------------------------------------------
# 4) The "Stipend for Single Degree" column (J) holds numeric-looking text
#    ("0") that must stay stored as text, matching the source data. Force
#    text storage with a leading quote-prefix (mirrors typing '0 in Excel).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 10).Value = "'0"
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
